$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Excel constant ---
$xlCenter = -4108

# --- Row 3: Vuong Thi Tiu (wife of 140001) -----------------------------
# Drop the row-level custom format (s="5" customFormat="1") entirely, then
# re-apply only the per-cell centering that should remain. The "Dinh" (C)
# marker cell is removed outright since this row is a wife, not a blood
# descendant.
$ws.Rows(3).ClearFormats()
$ws.Range("C3").Clear()
$ws.Range("A3").HorizontalAlignment = $xlCenter
$ws.Range("B3").HorizontalAlignment = $xlCenter
$ws.Range("J3").HorizontalAlignment = $xlCenter
$ws.Range("E3").HorizontalAlignment = $xlCenter
$ws.Range("F3").HorizontalAlignment = $xlCenter

# --- Row 4: Duong Thi Tan ------------------------------------------------
$ws.Rows(4).ClearFormats()
$ws.Range("C4").Clear()
$ws.Range("A4").HorizontalAlignment = $xlCenter
$ws.Range("B4").HorizontalAlignment = $xlCenter
$ws.Range("H4").HorizontalAlignment = $xlCenter
$ws.Range("E4").HorizontalAlignment = $xlCenter
$ws.Range("F4").HorizontalAlignment = $xlCenter
# new "mother ID" reference for this child of 140001v2
$ws.Range("I4").Value = "140001v2"
$ws.Range("I4").HorizontalAlignment = $xlCenter

# --- Row 5: Duong Trung Chinh --------------------------------------------
$ws.Range("H5").HorizontalAlignment = $xlCenter
$ws.Range("E5").HorizontalAlignment = $xlCenter
$ws.Range("I5").Value = "140001v2"
$ws.Range("I5").HorizontalAlignment = $xlCenter

# --- Row 6: empty placeholder "ID cha" cell removed ----------------------
$ws.Range("H6").Clear()
$ws.Range("E6").HorizontalAlignment = $xlCenter

# --- Row 7: Duong Thi Chinh -----------------------------------------------
$ws.Range("H7").HorizontalAlignment = $xlCenter
$ws.Range("E7").HorizontalAlignment = $xlCenter
$ws.Range("I7").Value = "140001v2"
$ws.Range("I7").HorizontalAlignment = $xlCenter

# --- Row 8: Duong Minh Tu ---------------------------------------------------
$ws.Range("H8").HorizontalAlignment = $xlCenter
$ws.Range("E8").HorizontalAlignment = $xlCenter
$ws.Range("I8").Value = "140001v2"
$ws.Range("I8").HorizontalAlignment = $xlCenter

# --- Header / summary rows: promote bold cells to bold+centered ----------
$ws.Range("E1").HorizontalAlignment = $xlCenter
$ws.Range("F1").HorizontalAlignment = $xlCenter
$ws.Range("E2").HorizontalAlignment = $xlCenter
$ws.Range("F2").HorizontalAlignment = $xlCenter

# --- Remaining "Nam sinh"/"Nam mat" (birth/death year) cells: center -----
$yearCells = @("E10","E11","E12","E13","E14","E15","E16","E17","E18","E19", `
               "E21","F21","E22","E23","E24","E25","E26","E27","E28","E29", `
               "E30","E31","E32")
foreach ($addr in $yearCells) {
    $ws.Range($addr).HorizontalAlignment = $xlCenter
}

# --- Selection cosmetic change --------------------------------------------
$ws.Range("G13").Select()
